$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range('D2')
$cD.Value = "'33.982.58"
$cD.Style = 'Normal'
$ws.Range('E2').Value = '  -0.32%  '
$cD = $ws.Range('D3')
$cD.Value = "'1.773.15"
$cD.Style = 'Normal'
$ws.Range('E3').Value = '  -2.26%  '
$cD = $ws.Range('D4')
$cD.Value = "'1.00"
$cD.Style = 'Normal'
$ws.Range('E4').Value = '  +0.28%  '
$cD = $ws.Range('D5')
$cD.Value = "'225.18"
$cD.Style = 'Normal'
$ws.Range('E5').Value = '  -1.25%  '
$cD = $ws.Range('D6')
$cD.Value = "'0.547"
$cD.Style = 'Normal'
$ws.Range('E6').Value = '  +0.53%  '
$cD = $ws.Range('D7')
$cD.Value = "'1.00"
$cD.Style = 'Normal'
$ws.Range('E7').Value = '  +0.17%  '
$cD = $ws.Range('D8')
$cD.Value = "'31.33"
$cD.Style = 'Normal'
$ws.Range('E8').Value = '  -0.08%  '
$cD = $ws.Range('D9')
$cD.Value = "'0.279"
$cD.Style = 'Normal'
$ws.Range('E9').Value = '  -0.48%  '
$cD = $ws.Range('D10')
$cD.Value = "'0.0655"
$cD.Style = 'Normal'
$ws.Range('E10').Value = '  -1.67%  '
$cD = $ws.Range('D11')
$cD.Value = "'0.0928"
$cD.Style = 'Normal'
$ws.Range('E11').Value = '  +0.01%  '
$cD = $ws.Range('D12')
$cD.Value = "'2.026.03"
$cD.Style = 'Normal'
$ws.Range('E12').Value = '  -2.29%  '
$cD = $ws.Range('D13')
$cD.Value = "'10.87"
$cD.Style = 'Normal'
$ws.Range('E13').Value = '  +6.54%  '
$cD = $ws.Range('D14')
$cD.Value = "'1.771.45"
$cD.Style = 'Normal'
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cD = $ws.Range('D15')
$cD.Value = "'0.621"
$cD.Style = 'Normal'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cD = $ws.Range('D16')
$cD.Value = "'33.977.06"
$cD.Style = 'Normal'
$ws.Range('E16').Value = '  -0.06%  '
$cD = $ws.Range('D17')
$cD.Value = "'4.19"
$cD.Style = 'Normal'
$ws.Range('E17').Value = '  -1.62%  '
$cD = $ws.Range('D18')
$cD.Value = "'68.55"
$cD.Style = 'Normal'
$ws.Range('E18').Value = '  -1.12%  '
$cD = $ws.Range('D19')
$cD.Value = "'251.53"
$cD.Style = 'Normal'
$ws.Range('E19').Value = '  -2.37%  '
$cD = $ws.Range('D20')
$cD.Value = "'0.0₃0735"
$cD.Style = 'Normal'
$ws.Range('E20').Value = '  -1.62%  '
$cD = $ws.Range('D21')
$cD.Value = "'1.00"
$cD.Style = 'Normal'
$ws.Range('E21').Value = '  +0.35%  '
$cD = $ws.Range('D22')
$cD.Value = "'10.30"
$cD.Style = 'Normal'
$ws.Range('E22').Value = '  -2.06%  '
$cD = $ws.Range('D23')
$cD.Value = "'4.18"
$cD.Style = 'Normal'
$ws.Range('E23').Value = '  -3.52%  '
$cD = $ws.Range('D24')
$cD.Value = "'2.12"
$cD.Style = 'Normal'
$ws.Range('E24').Value = '  -3.37%  '
$cD = $ws.Range('D25')
$cD.Value = "'156.45"
$cD.Style = 'Normal'
$ws.Range('E25').Value = '  -1.17%  '
$cD = $ws.Range('D26')
$cD.Value = "'16.34"
$cD.Style = 'Normal'
$ws.Range('E26').Value = '  -1.38%  '
$cD = $ws.Range('D27')
$cD.Value = "'6.98"
$cD.Style = 'Normal'
$ws.Range('E27').Value = '  -1.98%  '
$cD = $ws.Range('D28')
$cD.Value = "'0.113"
$cD.Style = 'Normal'
$ws.Range('E28').Value = '  -1.42%  '
$cD = $ws.Range('D29')
$cD.Value = "'1.00"
$cD.Style = 'Normal'
$ws.Range('E29').Value = '  +0.22%  '
$cD = $ws.Range('D30')
$cD.Value = "'3.74"
$cD.Style = 'Normal'
$ws.Range('E30').Value = '  -3.21%  '
$cD = $ws.Range('D31')
$cD.Value = "'0.0509"
$cD.Style = 'Normal'
$ws.Range('E31').Value = '  -0.43%  '
$cD = $ws.Range('D32')
$cD.Value = "'1.19"
$cD.Style = 'Normal'
$ws.Range('E32').Value = '  -0.62%  '
$cD = $ws.Range('D33')
$cD.Value = "'3.56"
$cD.Style = 'Normal'
$ws.Range('E33').Value = '  +1.47%  '
$cD = $ws.Range('D34')
$cD.Value = "'1.84"
$cD.Style = 'Normal'
$ws.Range('E34').Value = '  +2.09%  '
$cD = $ws.Range('D35')
$cD.Value = "'1.444.45"
$cD.Style = 'Normal'
$ws.Range('E35').Value = '  -6.09%  '
$cD = $ws.Range('D36')
$cD.Value = "'1.05"
$cD.Style = 'Normal'
$ws.Range('E36').Value = '  -2.05%  '
$cD = $ws.Range('D37')
$cD.Value = "'0.624"
$cD.Style = 'Normal'
$ws.Range('E37').Value = '  +0.18%  '
$cD = $ws.Range('D38')
$cD.Value = "'0.0186"
$cD.Style = 'Normal'
$ws.Range('E38').Value = '  -0.58%  '
$cD = $ws.Range('D39')
$cD.Value = "'2.83"
$cD.Style = 'Normal'
$ws.Range('E39').Value = '  +0.89%  '
$cD = $ws.Range('D40')
$cD.Value = "'82.48"
$cD.Style = 'Normal'
$ws.Range('E40').Value = '  -2.16%  '
$cD = $ws.Range('D41')
$cD.Value = "'2.34"
$cD.Style = 'Normal'
$ws.Range('E41').Value = '  +0.27%  '
$cD = $ws.Range('D42')
$cD.Value = "'0.885"
$cD.Style = 'Normal'
$ws.Range('E42').Value = '  -2.47%  '
$cD = $ws.Range('D43')
$cD.Value = "'2.03"
$cD.Style = 'Normal'
$ws.Range('E43').Value = '  -4.80%  '
$cD = $ws.Range('D44')
$cD.Value = "'0.0506"
$cD.Style = 'Normal'
$ws.Range('E44').Value = '  -2.61%  '
$cD = $ws.Range('D45')
$cD.Value = "'1.05"
$cD.Style = 'Normal'
$ws.Range('E45').Value = '  -1.83%  '
$cD = $ws.Range('D46')
$cD.Value = "'1.927.72"
$cD.Style = 'Normal'
$ws.Range('E46').Value = '  -2.03%  '
$cD = $ws.Range('D47')
$cD.Value = "'5.74"
$cD.Style = 'Normal'
$ws.Range('E47').Value = '  +0.19%  '
$cD = $ws.Range('D48')
$cD.Value = "'1.00"
$cD.Style = 'Normal'
$ws.Range('E48').Value = '  +0.30%  '
$cD = $ws.Range('D49')
$cD.Value = "'11.82"
$cD.Style = 'Normal'
$ws.Range('E49').Value = '  +2.12%  '
$cD = $ws.Range('D50')
$cD.Value = "'97.23"
$cD.Style = 'Normal'
$ws.Range('E50').Value = '  +2.76%  '
$cD = $ws.Range('D51')
$cD.Value = "'49.33"
$cD.Style = 'Normal'
$ws.Range('E51').Value = '  -5.91%  '
